$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("class_schedule")
$ws.Activate()

# Row 6 / Column C ("CI: Experiments" do-before-class reading):
# fix stray leading backtick typo ("`Internal" -> "Internal").
# Leading "'" forces literal text (keeps/produces the quotePrefix style flag,
# since the text still starts with "-").
$c6 = "'- Internal versus External Validity`n- ``Limitations of Average Treatment Effects <limitations_of_ATE.ipynb>``_`n- ``Experiments in Advertising <https://overcast.fm/+QLduPjO1I>``_"
$ws.Range("C6").Value = $c6

# Row 7 (CI: Natural Experiments) In-Class-Exercise column: add "- SUTVA" line.
$b7 = "'- SUTVA`n- Compliance / ITT / Etc."
$ws.Range("B7").Value = $b7
$ws.Range("B7").WrapText = $true

# Row 7 Do-Before-Class column: add Imbens & Rubin SUTVA reading + notation note.
$c7 = "'- Imbens and Rubin (CI), Section 1.6 (SUTVA, p. 10-13)`n- Angrist and Pischke (MM), Chapter 3 (pp 98-146)`n(Note that Imben & Rubin potential outcomes notation is a little different -- just skip notational parts if needed)"
$ws.Range("C7").Value = $c7

# Row 7 grew from a one-line row to a three-line row now that B7/C7 wrap;
# match the taller row height used once the new text is in place.
$ws.Rows.Item(7).RowHeight = 85

# Selection moved from C4 to C8 in the saved view.
$ws.Range("C8").Select()
